$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("5840560 - Marco Antonio Carvalho Pereira")) {
        $p.Range.Delete()
        break
    }
}
